$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header fields ---
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 25.10.2024"

# --- Row 6 ---
$ws.Range("B6").Value = "27.10."
$ws.Range("C6").Value = "28.10."
$ws.Range("D6").Value = "KARTENZ./27.10 LIDL RO"
$ws.Range("E6").Value = "141,16-"

# --- Row 7 ---
$ws.Range("B7").Value = "30.10."
$ws.Range("C7").Value = "31.10."
$ws.Range("D7").Value = "MCDONALDS Beilngries"
$ws.Range("E7").Value = "13,16-"

# --- Row 8 ---
$ws.Range("B8").Value = "02.11."
$ws.Range("C8").Value = "03.11."
$ws.Range("D8").Value = "PAYPAL LERBJK"
$ws.Range("E8").Value = "90,12-"

# --- Row 9 ---
$ws.Range("B9").Value = "04.11."
$ws.Range("C9").Value = "05.11."
$ws.Range("D9").Value = "PAYPAL MZEIDS"
$ws.Range("E9").Value = "52,00-"

# --- Row 10 ---
$ws.Range("B10").Value = "07.11."
$ws.Range("C10").Value = "08.11."
$ws.Range("D10").Value = "BEITRAG Allianz SE K-174949"
$ws.Range("E10").Value = "57,02-"

# --- Row 11: transaction removed -> clear contents ---
$ws.Range("B11:D11").ClearContents()
$ws.Range("E11").Value = ""
$ws.Range("E11").WrapText = $true
$ws.Range("E11").VerticalAlignment = -4108

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 09.11.2024"
$ws.Range("E12").Value = "353,46-"

# --- Next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 14.11.2024"
